# Apply results of new GMM function: fill in newly-computed correlation
# values in the adjacency matrix (previously placeholder zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "L3"  = -0.195520137061439
    "H6"  = -0.180854997919969
    "L6"  = 0.174918337523784
    "E9"  = -0.180854997919969
    "M12" = 0.192153395150204
    "U12" = 0.167664681960564
    "B13" = -0.195520137061439
    "E13" = 0.174918337523784
    "K14" = 0.192153395150204
    "S19" = -0.162793178746392
    "R20" = -0.162793178746392
    "K22" = 0.167664681960564
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
